$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.177.73"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "1.573.27"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.29"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("E6").Value = "  -1.76%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.21"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0589"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0868"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").Value = "1.797.19"
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("D13").Value = "1.574.61"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.518"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.24%  "
$ws.Range("D16").Value = "27.194.50"
$ws.Range("E16").Value = "  -1.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.22"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "214.13"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.35"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("D20").Value = "0.0₃0685"
$ws.Range("E20").Value = "  -0.90%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.38"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.89%  "
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.66"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("E26").Value = "  -3.65%  "
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  -1.26%  "
$ws.Range("E30").Value = "  -3.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0463"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.74%  "
$ws.Range("E32").Value = "  -1.24%  "
$ws.Range("D33").Value = "1.407.71"
$ws.Range("E33").Value = "  +3.30%  "
$ws.Range("E34").Value = "  -1.27%  "
$ws.Range("E35").Value = "  +1.89%  "
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.939"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.08%  "
$ws.Range("E38").Value = "  -1.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.816"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.515"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.995"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.82"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.61%  "
$ws.Range("E44").Value = "  +2.21%  "
$ws.Range("E45").Value = "  +1.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.72"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("D47").Value = "1.709.51"
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.86"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.60%  "
$ws.Range("D49").Value = "0.0₇0975"
$ws.Range("E49").Value = "  -2.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0953"
$ws.Range("D50").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0494"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.17%  "
